# Add question_code and result_type columns to the Lookups sheet,
# mirroring the existing Gender/male/female "cuts" columns (F:G) into
# two new column pairs (H:I and J:K), and widen the cuts_head named
# range so it covers the new header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# --- Row 1 (headers): add question_code / result_type headers ---
$ws.Range("H1").Value = "question_code"
$ws.Range("J1").Value = "result_type"

# --- Row 2: mirror the existing Gender/male row (F2:G2) into H2:I2 and J2:K2 ---
$ws.Range("H2").Value = $ws.Range("F2").Value2
$ws.Range("I2").Value = $ws.Range("G2").Value2
$ws.Range("J2").Value = $ws.Range("F2").Value2
$ws.Range("K2").Value = $ws.Range("G2").Value2

# --- Row 3: mirror the existing female row (F3:G3) into H3:I3 and J3:K3 ---
$ws.Range("H3").Value = $ws.Range("F3").Value2
$ws.Range("I3").Value = $ws.Range("G3").Value2
$ws.Range("J3").Value = $ws.Range("F3").Value2
$ws.Range("K3").Value = $ws.Range("G3").Value2

# --- Widen the "cuts_head" defined name so it spans F1:K1 instead of F1:G1 ---
$wb.Names.Item("cuts_head").RefersTo = "='Lookups'!`$F`$1:`$K`$1"
